$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='295.59'},
    @{Cell='E2'; Value='1.45%'},
    @{Cell='F2'; Value='15-2-2023'},
    @{Cell='G2'; Value='0'},
    @{Cell='D3'; Value='42.17'},
    @{Cell='E3'; Value='3.25%'},
    @{Cell='F3'; Value='15-2-2023'},
    @{Cell='G3'; Value='0'},
    @{Cell='D4'; Value='5.022'},
    @{Cell='E4'; Value='-0.40%'},
    @{Cell='F4'; Value='15-2-2023'},
    @{Cell='G4'; Value='0'},
    @{Cell='D5'; Value='0.07556'},
    @{Cell='E5'; Value='2.44%'},
    @{Cell='F5'; Value='15-2-2023'},
    @{Cell='G5'; Value='0'},
    @{Cell='B6'; Value='GateToken'},
    @{Cell='C6'; Value='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'},
    @{Cell='D6'; Value='4.388'},
    @{Cell='E6'; Value='2.36%'},
    @{Cell='F6'; Value='15-2-2023'},
    @{Cell='G6'; Value='0'},
    @{Cell='B7'; Value='FTXToken'},
    @{Cell='C7'; Value='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'},
    @{Cell='D7'; Value='1.593'},
    @{Cell='E7'; Value='1.24%'},
    @{Cell='F7'; Value='15-2-2023'},
    @{Cell='G7'; Value='0'},
    @{Cell='B8'; Value='MXToken'},
    @{Cell='C8'; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Cell='D8'; Value='0.9282'},
    @{Cell='E8'; Value='0.35%'},
    @{Cell='F8'; Value='15-2-2023'},
    @{Cell='G8'; Value='0'},
    @{Cell='B9'; Value='BTSEToken'},
    @{Cell='C9'; Value='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'},
    @{Cell='D9'; Value='2.410'},
    @{Cell='E9'; Value='1.18%'},
    @{Cell='F9'; Value='15-2-2023'},
    @{Cell='G9'; Value='0'},
    @{Cell='B10'; Value='LiechtensteinCryptoassetsExchange'},
    @{Cell='C10'; Value='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'},
    @{Cell='D10'; Value='0.1200'},
    @{Cell='E10'; Value='4.30%'},
    @{Cell='F10'; Value='15-2-2023'},
    @{Cell='G10'; Value='0'},
    @{Cell='B11'; Value='WazirX'},
    @{Cell='C11'; Value='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'},
    @{Cell='D11'; Value='0.1836'},
    @{Cell='E11'; Value='5.30%'},
    @{Cell='F11'; Value='15-2-2023'},
    @{Cell='G11'; Value='0'},
    @{Cell='B12'; Value='MandalaExchangeToken'},
    @{Cell='C12'; Value='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'},
    @{Cell='D12'; Value='0.09007'},
    @{Cell='E12'; Value='3.52%'},
    @{Cell='F12'; Value='15-2-2023'},
    @{Cell='G12'; Value='0'},
    @{Cell='B13'; Value='BitrueCoin'},
    @{Cell='C13'; Value='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'},
    @{Cell='D13'; Value='0.04083'},
    @{Cell='E13'; Value='-2.31%'},
    @{Cell='F13'; Value='15-2-2023'},
    @{Cell='G13'; Value='0'},
    @{Cell='B14'; Value='BitMartToken'},
    @{Cell='C14'; Value='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'},
    @{Cell='D14'; Value='0.1051'},
    @{Cell='E14'; Value='-0.46%'},
    @{Cell='F14'; Value='15-2-2023'},
    @{Cell='G14'; Value='0'},
    @{Cell='B15'; Value='BitForexToken'},
    @{Cell='C15'; Value='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'},
    @{Cell='D15'; Value='0.001293'},
    @{Cell='E15'; Value='1.75%'},
    @{Cell='F15'; Value='15-2-2023'},
    @{Cell='G15'; Value='0'},
    @{Cell='B16'; Value='TigerCash'},
    @{Cell='C16'; Value='https://coinranking.com/coin/6hIn06L2+tigercash-tch'},
    @{Cell='D16'; Value='0.005876'},
    @{Cell='E16'; Value='-1.53%'},
    @{Cell='F16'; Value='15-2-2023'},
    @{Cell='G16'; Value='0'},
    @{Cell='B17'; Value='LEO'},
    @{Cell='C17'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'},
    @{Cell='D17'; Value='3.365'},
    @{Cell='E17'; Value='-0.83%'},
    @{Cell='F17'; Value='15-2-2023'},
    @{Cell='G17'; Value='0'},
    @{Cell='E18'; Value='1.07%'},
    @{Cell='F18'; Value='15-2-2023'},
    @{Cell='G18'; Value='0'},
    @{Cell='D19'; Value='7.833'},
    @{Cell='E19'; Value='1.29%'},
    @{Cell='F19'; Value='15-2-2023'},
    @{Cell='G19'; Value='0'},
    @{Cell='D20'; Value='0.1422'},
    @{Cell='E20'; Value='3.18%'},
    @{Cell='F20'; Value='15-2-2023'},
    @{Cell='G20'; Value='0'},
    @{Cell='D21'; Value='0.3004'},
    @{Cell='E21'; Value='4.18%'},
    @{Cell='F21'; Value='15-2-2023'},
    @{Cell='G21'; Value='0'},
    @{Cell='D22'; Value='0.04057'},
    @{Cell='E22'; Value='4.33%'},
    @{Cell='F22'; Value='15-2-2023'},
    @{Cell='G22'; Value='0'},
    @{Cell='D23'; Value='0.001266'},
    @{Cell='E23'; Value='0.63%'},
    @{Cell='F23'; Value='15-2-2023'},
    @{Cell='G23'; Value='0'},
    @{Cell='D24'; Value='0.004105'},
    @{Cell='E24'; Value='4.96%'},
    @{Cell='F24'; Value='15-2-2023'},
    @{Cell='G24'; Value='0'},
    @{Cell='E25'; Value='-3.92%'},
    @{Cell='F25'; Value='15-2-2023'},
    @{Cell='G25'; Value='0'},
    @{Cell='E26'; Value='-0.04%'},
    @{Cell='F26'; Value='15-2-2023'},
    @{Cell='G26'; Value='0'},
    @{Cell='F27'; Value='15-2-2023'},
    @{Cell='G27'; Value='0'},
    @{Cell='F28'; Value='15-2-2023'},
    @{Cell='G28'; Value='0'},
    @{Cell='F29'; Value='15-2-2023'},
    @{Cell='G29'; Value='0'},
    @{Cell='F30'; Value='15-2-2023'},
    @{Cell='G30'; Value='0'},
    @{Cell='F31'; Value='15-2-2023'},
    @{Cell='G31'; Value='0'},
    @{Cell='F32'; Value='15-2-2023'},
    @{Cell='G32'; Value='0'},
    @{Cell='F33'; Value='15-2-2023'},
    @{Cell='G33'; Value='0'},
    @{Cell='F34'; Value='15-2-2023'},
    @{Cell='G34'; Value='0'},
    @{Cell='F35'; Value='15-2-2023'},
    @{Cell='G35'; Value='0'},
    @{Cell='F36'; Value='15-2-2023'},
    @{Cell='G36'; Value='0'},
    @{Cell='F37'; Value='15-2-2023'},
    @{Cell='G37'; Value='0'},
    @{Cell='D38'; Value='0.02410'},
    @{Cell='E38'; Value='3.12%'},
    @{Cell='F38'; Value='15-2-2023'},
    @{Cell='G38'; Value='0'},
    @{Cell='D39'; Value='0.05207'},
    @{Cell='E39'; Value='3.74%'},
    @{Cell='F39'; Value='15-2-2023'},
    @{Cell='G39'; Value='0'},
    @{Cell='D40'; Value='0.005996'},
    @{Cell='E40'; Value='0.66%'},
    @{Cell='F40'; Value='15-2-2023'},
    @{Cell='G40'; Value='0'},
    @{Cell='D41'; Value='0.007774'},
    @{Cell='E41'; Value='1.18%'},
    @{Cell='F41'; Value='15-2-2023'},
    @{Cell='G41'; Value='0'},
    @{Cell='D42'; Value='0.1330'},
    @{Cell='E42'; Value='3.27%'},
    @{Cell='F42'; Value='15-2-2023'},
    @{Cell='G42'; Value='0'},
    @{Cell='D43'; Value='0.007564'},
    @{Cell='E43'; Value='2.94%'},
    @{Cell='F43'; Value='15-2-2023'},
    @{Cell='G43'; Value='0'},
    @{Cell='D44'; Value='0.007284'},
    @{Cell='E44'; Value='2.83%'},
    @{Cell='F44'; Value='15-2-2023'},
    @{Cell='G44'; Value='0'},
    @{Cell='D45'; Value='0.2939'},
    @{Cell='E45'; Value='-6.80%'},
    @{Cell='F45'; Value='15-2-2023'},
    @{Cell='G45'; Value='0'},
    @{Cell='D46'; Value='0.00006742'},
    @{Cell='E46'; Value='5.07%'},
    @{Cell='F46'; Value='15-2-2023'},
    @{Cell='G46'; Value='0'},
    @{Cell='D47'; Value='0.00000000751'},
    @{Cell='E47'; Value='0.00%'},
    @{Cell='F47'; Value='15-2-2023'},
    @{Cell='G47'; Value='0'},
    @{Cell='D48'; Value='0.04555'},
    @{Cell='E48'; Value='164.97%'},
    @{Cell='F48'; Value='15-2-2023'},
    @{Cell='G48'; Value='0'},
    @{Cell='D49'; Value='0.004208'},
    @{Cell='E49'; Value='0.14%'},
    @{Cell='F49'; Value='15-2-2023'},
    @{Cell='G49'; Value='0'},
    @{Cell='D50'; Value='0.00002102'},
    @{Cell='E50'; Value='0.00%'},
    @{Cell='F50'; Value='15-2-2023'},
    @{Cell='G50'; Value='0'},
    @{Cell='D51'; Value='0.0002002'},
    @{Cell='E51'; Value='0.00%'},
    @{Cell='F51'; Value='15-2-2023'},
    @{Cell='G51'; Value='0'}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}

Write-Host "Applied" $updates.Count "cell updates"
